$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing value in B42
$ws.Range("B42").Value = 28.31

# Populate new data for row 50 (year 2029)
$row50 = @{
    "B50" = 28.05
    "C50" = 29.14
    "D50" = 29.4
    "E50" = 29.5
    "F50" = 28.74
    "G50" = 27.54
    "H50" = 26.04
    "I50" = 25.49
    "J50" = 26.73
    "K50" = 26.56
    "L50" = 28.11
    "M50" = 28.4
}
foreach ($addr in $row50.Keys) {
    $ws.Range($addr).Value = $row50[$addr]
}

# Populate new data for row 51 (year 2030)
$row51 = @{
    "B51" = 27.97
    "C51" = 28.66
    "D51" = 29.38
    "E51" = 29.05
    "F51" = 28.81
    "G51" = 27.32
    "H51" = 26.01
    "I51" = 25.41
    "J51" = 26.14
    "K51" = 27.87
    "L51" = 27.34
    "M51" = 28.1
}
foreach ($addr in $row51.Keys) {
    $ws.Range($addr).Value = $row51[$addr]
}
